# Regenerate merged AHB files
# 1) Rename the header labels: "_old" -> "_FV2404", "_new" -> "_FV2410"
# 2) Freeze the header row (pane split below row 1)
# 3) Turn the used range into an Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# Columns A:J (1-10) hold the "_old" -> "_FV2404" headers
for ($i = 0; $i -lt $fv2404Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2404Headers[$i]
}

# Column K (11) is "diff" - unchanged.
# Columns L:U (12-21) hold the "_new" -> "_FV2410" headers
for ($i = 0; $i -lt $fv2410Headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2410Headers[$i]
}

# Freeze panes below the header row
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)

# Turn the used range (header + 57 data rows, 21 columns) into an Excel Table
$tableRange = $ws.Range("A1:U58")
$listObject = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObject.Name = "Table1"
